# #272 Ajout d'un scenario de recherche de l'offre d'un professionnel avec un ID Nat PS (#340)
# Applies:
#   1. Metadata!B8 "Date" value bump.
#   2. Elements sheet: swap the two "Mapping" columns (AK <-> AL), content + width,
#      so the "Specification metier" mapping now comes before "RIM Mapping".

$wb = $excel.ActiveWorkbook

# --- 1. Metadata: bump the Date value ---------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements: swap columns AK (37) and AL (38) ---------------------
$elements = $wb.Worksheets.Item("Elements")

# Only rows whose AK/AL pair actually differ need touching (rows 2 and 4
# hold an empty string in both columns already, so skip them to avoid
# needlessly rewriting those cells).
$rowsToSwap = @(1, 3, 5, 6)
foreach ($r in $rowsToSwap) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths to match (AK becomes the wide column, AL the narrow one).
# Target stored widths: AK -> 87.7890625, AL -> 24.98046875 (the previous AK width).
# NB: ColumnWidth is expressed in character units and gets re-quantised by the
# host to whole pixels on write, so the literals below are chosen to land on
# the closest achievable stored width to the target.
$elements.Columns.Item(37).ColumnWidth = 87.0
$elements.Columns.Item(38).ColumnWidth = 24.15
